# Update FFXIV leve-profit calculation sheets (currentAveragePrice* / LevePrice* / LeveProfit*)
# with refreshed market-board figures from the scheduled price-data runner.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 138: All-night Crafting
$ws.Range("H138").Value = 2981.75
$ws.Range("I138").Value = 2865.6667
$ws.Range("J138").Value = 2985.8472
$ws.Range("K138").Value = 8597.000100000001
$ws.Range("L138").Value = 8957.5416
$ws.Range("M138").Value = -3457.000100000001
$ws.Range("N138").Value = -19237.5416

$ws = $wb.Worksheets.Item("ARM")
# Row 2: Ain't Got No Ingots
$ws.Range("H2").Value = 333.57144
$ws.Range("I2").Value = 222.66667
$ws.Range("K2").Value = 222.66667
$ws.Range("M2").Value = -109.66667

# Row 92: Mail It In
$ws.Range("H92").Value = 104995
$ws.Range("J92").Value = 104995
$ws.Range("L92").Value = 104995
$ws.Range("N92").Value = -109987

# Row 114: A New Regular
$ws.Range("H114").Value = 50333
$ws.Range("J114").Value = 50333
$ws.Range("L114").Value = 50333
$ws.Range("N114").Value = -59011

# Row 116: No Scope
$ws.Range("H116").Value = 333.57144
$ws.Range("I116").Value = 222.66667
$ws.Range("K116").Value = 222.66667
$ws.Range("M116").Value = 2071.33333

# Row 122: Haste for High Durium
$ws.Range("H122").Value = 2524.5483
$ws.Range("I122").Value = 1806.7084
$ws.Range("K122").Value = 5420.1252
$ws.Range("M122").Value = -2970.1252

# Row 125: The Incomplete Costume
$ws.Range("H125").Value = 51244.75
$ws.Range("J125").Value = 51244.75
$ws.Range("L125").Value = 51244.75
$ws.Range("N125").Value = -61084.75

$ws = $wb.Worksheets.Item("BSM")
# Row 3: Hells Bells
$ws.Range("H3").Value = 333.57144
$ws.Range("I3").Value = 222.66667
$ws.Range("K3").Value = 222.66667
$ws.Range("M3").Value = -108.66667

# Row 7: Thank You for Your Business
$ws.Range("H7").Value = 13337728
$ws.Range("I7").Value = 3184
$ws.Range("K7").Value = 3184
$ws.Range("M7").Value = -3071

# Row 82: Spirituality Inspector
$ws.Range("H82").Value = 30631.2
$ws.Range("I82").Value = 6289
$ws.Range("J82").Value = 128000
$ws.Range("K82").Value = 6289
$ws.Range("L82").Value = 128000
$ws.Range("M82").Value = -5906
$ws.Range("N82").Value = -128766

# Row 85: The Clamor for Hammers (L)
$ws.Range("H85").Value = 30631.2
$ws.Range("I85").Value = 6289
$ws.Range("J85").Value = 128000
$ws.Range("K85").Value = 6289
$ws.Range("L85").Value = 128000
$ws.Range("M85").Value = -4963
$ws.Range("N85").Value = -130652

# Row 99: Meddle in Metal
$ws.Range("H99").Value = 6959.7827
$ws.Range("I99").Value = 9174.154
$ws.Range("K99").Value = 9174.154
$ws.Range("M99").Value = -7676.154

# Row 109: Here Comes the Hammer
$ws.Range("H109").Value = 103914.5
$ws.Range("J109").Value = 103914.5
$ws.Range("L109").Value = 103914.5
$ws.Range("N109").Value = -106688.5

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Range("H31").Value = 1310843.1
$ws.Range("I31").Value = 26421.666
$ws.Range("K31").Value = 26421.666
$ws.Range("M31").Value = -26126.666

# Row 34: Armoires of the Rich and Famous
$ws.Range("H34").Value = 1310843.1
$ws.Range("I34").Value = 26421.666
$ws.Range("K34").Value = 26421.666
$ws.Range("M34").Value = -26219.666

# Row 99: O Pine
$ws.Range("H99").Value = 2953.3333
$ws.Range("J99").Value = 3303.6667
$ws.Range("L99").Value = 3303.6667
$ws.Range("N99").Value = -6299.6667

# Row 126: A Better Conductor
$ws.Range("H126").Value = 2953.3333
$ws.Range("J126").Value = 3303.6667
$ws.Range("L126").Value = 9911.000100000001
$ws.Range("N126").Value = -14851.0001

$ws = $wb.Worksheets.Item("CUL")
# Row 33: Cooking with Gas
$ws.Range("H33").Value = 525
$ws.Range("I33").Value = 50
$ws.Range("K33").Value = 300
$ws.Range("M33").Value = -17

# Row 98: Sweet Kiss of Death
$ws.Range("H98").Value = 236.5
$ws.Range("I98").Value = 236.5
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 709.5
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = $null
$ws.Range("N98").Value = 788.5

$ws = $wb.Worksheets.Item("GSM")
# Row 80: Needs More Prayerbell
$ws.Range("H80").Value = 4735.75
$ws.Range("I80").Value = 4131.875
$ws.Range("J80").Value = 5943.5
$ws.Range("K80").Value = 4131.875
$ws.Range("L80").Value = 5943.5
$ws.Range("M80").Value = -3133.875
$ws.Range("N80").Value = -7939.5

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Range("H83").Value = 4735.75
$ws.Range("I83").Value = 4131.875
$ws.Range("J83").Value = 5943.5
$ws.Range("K83").Value = 20659.375
$ws.Range("L83").Value = 29717.5
$ws.Range("M83").Value = -15667.375
$ws.Range("N83").Value = -39701.5

# Row 126: Gold Rush Order
$ws.Range("H126").Value = 4421.364
$ws.Range("I126").Value = 3939.1667
$ws.Range("K126").Value = 11817.5001
$ws.Range("M126").Value = -9347.500100000001

# Row 136: Shiny and Good
$ws.Range("H136").Value = 6262.773
$ws.Range("J136").Value = 6262.773
$ws.Range("L136").Value = 18788.319
$ws.Range("N136").Value = -23888.319

$ws = $wb.Worksheets.Item("LTW")
# Row 7: Tan Before the Ban
$ws.Range("H7").Value = 588603.8
$ws.Range("I7").Value = 10347.9
$ws.Range("J7").Value = 1231110.4
$ws.Range("K7").Value = 10347.9
$ws.Range("L7").Value = 1231110.4
$ws.Range("M7").Value = -10235.9
$ws.Range("N7").Value = -1231334.4

# Row 16: Saddle Sore
$ws.Range("H16").Value = 912.125
$ws.Range("I16").Value = 622.38464
$ws.Range("K16").Value = 622.38464
$ws.Range("M16").Value = -452.38464

# Row 100: Tiger in the Sack
$ws.Range("H100").Value = 2799.3
$ws.Range("I100").Value = 1149.25
$ws.Range("K100").Value = 1149.25
$ws.Range("M100").Value = -608.25

# Row 109: Band Substances
$ws.Range("H109").Value = 80500
$ws.Range("J109").Value = 80500
$ws.Range("L109").Value = 80500
$ws.Range("N109").Value = -83274

# Row 126: Battered Books
$ws.Range("H126").Value = 588603.8
$ws.Range("I126").Value = 10347.9
$ws.Range("J126").Value = 1231110.4
$ws.Range("K126").Value = 31043.7
$ws.Range("L126").Value = 3693331.2
$ws.Range("M126").Value = -28573.7
$ws.Range("N126").Value = -3698271.2

$ws = $wb.Worksheets.Item("WVR")
# Row 2: The Unmentionables
$ws.Range("H2").Value = 628029.3
$ws.Range("I2").Value = 711333.2
$ws.Range("J2").Value = 3250
$ws.Range("K2").Value = 711333.2
$ws.Range("L2").Value = 3250
$ws.Range("M2").Value = -711221.2
$ws.Range("N2").Value = -3474

# Row 26: New Shoes, New Me
$ws.Range("H26").Value = 10000
$ws.Range("I26").Value = 0
$ws.Range("K26").Value = 0
$ws.Range("M26").Value = $null

# Row 53: I'll Swap You
$ws.Range("H53").Value = 41250
$ws.Range("I53").Value = 40000
$ws.Range("J53").Value = 42500
$ws.Range("K53").Value = 40000
$ws.Range("L53").Value = 42500
$ws.Range("M53").Value = -39393
$ws.Range("N53").Value = -43714

# Row 81: Where the Dragonflies, the Net Catches
$ws.Range("H81").Value = 7161.6787
$ws.Range("I81").Value = 3668.8125
$ws.Range("J81").Value = 11818.833
$ws.Range("K81").Value = 7337.625
$ws.Range("L81").Value = 23637.666
$ws.Range("M81").Value = -6276.625
$ws.Range("N81").Value = -25759.666

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Range("H84").Value = 7161.6787
$ws.Range("I84").Value = 3668.8125
$ws.Range("J84").Value = 11818.833
$ws.Range("K84").Value = 36688.125
$ws.Range("L84").Value = 118188.33
$ws.Range("M84").Value = -31384.125
$ws.Range("N84").Value = -128796.33

# Row 109: Turban in Training
$ws.Range("H109").Value = 106989.5
$ws.Range("J109").Value = 106989.5
$ws.Range("L109").Value = 106989.5
$ws.Range("N109").Value = -109763.5
